$wb = $excel.ActiveWorkbook

# --- Append a new log row to the "Logs" sheet ---
$logs = $wb.Worksheets.Item("Logs")

$newRow = 23
$logs.Cells.Item($newRow, 1).Value = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-13 22:10:01"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend the conditional formatting ranges so they cover the new row ---
$cfCols = @("D","G","H","I","J")
foreach ($col in $cfCols) {
  $target = $logs.Range($col + "2:" + $col + $newRow)
  $fc = $target.FormatConditions.Item(1)
  $fc.ModifyAppliesToRange($target)
}

# --- Update the Dashboard summary count for this category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 22
